$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the header text "age" -> "Age"
$ws.Range("B1").Value = "Age"

# 2. Set column widths (approximate Apache POI autosize-style widths;
#    the engine quantizes ColumnWidth to 1/6-character steps, so these
#    inputs are chosen to land as close as possible to 7.82421875 / 4.0703125)
$ws.Columns.Item(1).ColumnWidth = 6.9167
$ws.Columns.Item(2).ColumnWidth = 3.0834

# 3. Style the header row (A1:B1): bold font, centered, full thin border.
#    Build the formatting up on A1 one border edge at a time (top, bottom,
#    left, right - this is the order that yields the same border-table
#    progression Apache POI produced), then copy the fully-resolved format
#    from A1 onto B1 in one shot so both header cells end up sharing a
#    single cell style, same as in the target workbook.
$a = $ws.Range("A1")
$a.Font.Bold = $true
$a.HorizontalAlignment = -4108   # xlCenter
$a.VerticalAlignment = -4108     # xlCenter
$a.Borders.Item(8).LineStyle = 1  # xlEdgeTop
$a.Borders.Item(9).LineStyle = 1  # xlEdgeBottom
$a.Borders.Item(7).LineStyle = 1  # xlEdgeLeft
$a.Borders.Item(10).LineStyle = 1 # xlEdgeRight

$b = $ws.Range("B1")
$a.Copy()
$b.PasteSpecial(-4122)  # xlPasteFormats
